# Add "Dead Code" error entry to the Errors table (row 20, ID=8, SemAnalyzer stage)
# and move the active selection to D20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors")

# Row 20 currently has Severity=Error, ID=8, Stage=SemAnalyzer, Description=<empty>.
# Update Severity to "Warning" and fill in the Description ("Dead code").
$ws.Range("A20").Value = "Warning"
$ws.Range("D20").Value = "Dead code"

# Move the selection cursor to D20 (matches the diff's sheetView selection change).
$ws.Range("D20").Select()
